# Update "Pin Mapping" worksheet: rename the LLBV3 Header column to
# "LLBV3 Header / Function" and fill in the function/purpose notes that
# were missing for several pins, matching the current state of the board.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header
$ws.Cells.Item(1, 5).Value = "LLBV3 Header / Function"

# New notes, entered in the same order the author typed them so brand new
# entries land in the shared-string table in the right sequence.
$ws.Cells.Item(2, 5).Value  = "MCP 2515 interrupt on received frames"
$ws.Cells.Item(8, 5).Value  = "E-stop jumper, also X3"
$ws.Cells.Item(19, 5).Value = "wheel hall switch header"
$ws.Cells.Item(21, 5).Value = "all SPI devices, SPI header"
$ws.Cells.Item(22, 5).Value = "all SPI devices, SPI header"
$ws.Cells.Item(23, 5).Value = "all SPI devices, SPI header"
$ws.Cells.Item(20, 5).Value = "SPI header (this pin tells the mega to be a slave)"
$ws.Cells.Item(36, 5).Value = "MCP2515 slave selection"
$ws.Cells.Item(37, 5).Value = "DAC slave selection"
$ws.Cells.Item(53, 5).Value = "on-board buzzer"
$ws.Cells.Item(55, 5).Value = "X3, no purpose assgined"

# Remaining cells reuse notes that already exist elsewhere in the sheet.
$ws.Cells.Item(3, 5).Value  = "USB Serial"
$ws.Cells.Item(4, 5).Value  = "USB Serial"
$ws.Cells.Item(6, 5).Value  = "X3, for power on board"
$ws.Cells.Item(7, 5).Value  = "X3, for power on board"
$ws.Cells.Item(18, 5).Value = "X3, for power on board"
$ws.Cells.Item(24, 5).Value = "X3, for power on board"
$ws.Cells.Item(25, 5).Value = "X3, for power on board"
$ws.Cells.Item(26, 5).Value = "X3, for power on board"
$ws.Cells.Item(27, 5).Value = "X3, for power on board"
$ws.Cells.Item(52, 5).Value = "on-board relay"
$ws.Cells.Item(54, 5).Value = "on-board relay"
$ws.Cells.Item(57, 5).Value = "X3, for power on board"
$ws.Cells.Item(64, 5).Value = "X3, for power on board"
$ws.Cells.Item(65, 5).Value = "X3, for power on board"
$ws.Cells.Item(77, 5).Value = "X3, for power on board"
$ws.Cells.Item(79, 5).Value = "X3, for power on board"
$ws.Cells.Item(17, 5).Value = "Steering header"
$ws.Cells.Item(91, 5).Value = "Steering header"
$ws.Cells.Item(92, 5).Value = "Steering header"
$ws.Cells.Item(95, 5).Value = "Steering header"
$ws.Cells.Item(96, 5).Value = "Steering header"

# Move the view back up to where the author was working.
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("E58").Select()
